$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend the "ano" / "domicilios" series from 2050 through 2060 (rows 43-52)
$newData = @(
    @(2051, 101899855.555484),
    @(2052, 102620207.683704),
    @(2053, 103340559.81192499),
    @(2054, 104060911.940146),
    @(2055, 104781264.06836601),
    @(2056, 105501616.196587),
    @(2057, 106221968.324808),
    @(2058, 106942320.45302799),
    @(2059, 107662672.581249),
    @(2060, 108383024.70947)
)

$row = 43
foreach ($pair in $newData) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}
